$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  19"
$ws.Range("C9").Value = "Report Covering the Week  5/8/2023  Through  5/14/2023"

# --- Fix up cells whose data type changes (text <-> number) by copying the
#     exact format (and, for str->str donors, content) from a stable donor cell
#     elsewhere in the sheet whose own value the edit leaves untouched. ---
$ws.Range("I14").Copy($ws.Range("C15"))
$ws.Range("I14").Copy($ws.Range("D16"))
$ws.Range("N14").Copy($ws.Range("E16"))
$ws.Range("C14").Copy($ws.Range("D18"))  # copy supplies final text value too
$ws.Range("E14").Copy($ws.Range("E18"))  # copy supplies final text value too
$ws.Range("I14").Copy($ws.Range("C26"))
$ws.Range("I14").Copy($ws.Range("C27"))
$ws.Range("I14").Copy($ws.Range("D27"))
$ws.Range("N14").Copy($ws.Range("E27"))
$ws.Range("N14").Copy($ws.Range("N28"))
$ws.Range("N14").Copy($ws.Range("N29"))
$ws.Range("D14").Copy($ws.Range("C30"))  # copy supplies final text value too

# --- Apply final values for all changed numeric cells (and the few cells that
#     become new text values not already supplied verbatim by the donor copy) ---
$ws.Range("C15").Value = 2
$ws.Range("F15").Value = 2
$ws.Range("I15").Value = 5
$ws.Range("M15").Value = 150
$ws.Range("N15").Value = 0
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = -40
$ws.Range("I16").Value = 23
$ws.Range("J16").Value = 34
$ws.Range("K16").Value = -32.35294117647
$ws.Range("L16").Value = 283.333333333333
$ws.Range("M16").Value = -23.333333333333
$ws.Range("N16").Value = -82.575757575757
$ws.Range("C17").Value = 2
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 11
$ws.Range("H17").Value = 120
$ws.Range("I17").Value = 37
$ws.Range("J17").Value = 26
$ws.Range("K17").Value = 42.307692307692
$ws.Range("L17").Value = 117.647058823529
$ws.Range("M17").Value = 76.190476190476
$ws.Range("N17").Value = -22.916666666666
$ws.Range("C18").Value = 10
$ws.Range("F18").Value = 27
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 170
$ws.Range("I18").Value = 129
$ws.Range("K18").Value = 21.698113207547
$ws.Range("L18").Value = 57.317073170731
$ws.Range("M18").Value = 40.217391304347
$ws.Range("N18").Value = -64.560439560439
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 66.666666666666
$ws.Range("F19").Value = 54
$ws.Range("G19").Value = 41
$ws.Range("H19").Value = 31.70731707317
$ws.Range("I19").Value = 273
$ws.Range("J19").Value = 199
$ws.Range("K19").Value = 37.185929648241
$ws.Range("L19").Value = 123.770491803279
$ws.Range("M19").Value = 92.25352112676
$ws.Range("N19").Value = 50
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 13
$ws.Range("H20").Value = 160
$ws.Range("I20").Value = 54
$ws.Range("J20").Value = 36
$ws.Range("K20").Value = 50
$ws.Range("L20").Value = 157.142857142857
$ws.Range("M20").Value = 10.204081632653
$ws.Range("N20").Value = -95.423728813559
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 12
$ws.Range("E21").Value = 175
$ws.Range("F21").Value = 110
$ws.Range("G21").Value = 66
$ws.Range("H21").Value = 66.666666666666
$ws.Range("I21").Value = 522
$ws.Range("J21").Value = 401
$ws.Range("K21").Value = 30.174563591022
$ws.Range("L21").Value = 110.483870967742
$ws.Range("M21").Value = 55.357142857142
$ws.Range("N21").Value = -72.698744769874
$ws.Range("C24").Value = 8
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = -46.666666666666
$ws.Range("F24").Value = 42
$ws.Range("G24").Value = 70
$ws.Range("H24").Value = -40
$ws.Range("I24").Value = 197
$ws.Range("J24").Value = 305
$ws.Range("K24").Value = -35.409836065573
$ws.Range("L24").Value = 2.072538860103
$ws.Range("M24").Value = 7.650273224043
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -75
$ws.Range("G25").Value = 16
$ws.Range("H25").Value = 12.5
$ws.Range("I25").Value = 77
$ws.Range("J25").Value = 79
$ws.Range("K25").Value = -2.53164556962
$ws.Range("L25").Value = 92.5
$ws.Range("M25").Value = 40
$ws.Range("C26").Value = 2
$ws.Range("F26").Value = 2
$ws.Range("I26").Value = 6
$ws.Range("L26").Value = 500
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 7
$ws.Range("J27").Value = 9
$ws.Range("K27").Value = -22.222222222222
$ws.Range("L27").Value = 40
$ws.Range("N28").Value = 0
$ws.Range("N29").Value = 0

# --- Column E width adjustment (closest achievable approximation of the bestFit
#     recompute Excel performed after the underlying data changed widths) ---
$ws.Columns("E:E").ColumnWidth = 8
